$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "43+28="
$t.Cell(1,2).Range.Text = "83-9="
$t.Cell(1,3).Range.Text = "19+12="
$t.Cell(1,4).Range.Text = "46+37="
$t.Cell(1,5).Range.Text = "37+25="
$t.Cell(2,1).Range.Text = "94-49="
$t.Cell(2,2).Range.Text = "25-19="
$t.Cell(2,3).Range.Text = "34+49="
$t.Cell(2,4).Range.Text = "56+28="
$t.Cell(2,5).Range.Text = "45+26="
$t.Cell(3,1).Range.Text = "39+53="
$t.Cell(3,2).Range.Text = "93-68="
$t.Cell(3,3).Range.Text = "55+16="
$t.Cell(3,4).Range.Text = "43-8="
$t.Cell(3,5).Range.Text = "22+9="
$t.Cell(4,1).Range.Text = "27+39="
$t.Cell(4,2).Range.Text = "70-53="
$t.Cell(4,3).Range.Text = "77+7="
$t.Cell(4,4).Range.Text = "7+36="
$t.Cell(4,5).Range.Text = "16+57="
$t.Cell(5,1).Range.Text = "36+9="
$t.Cell(5,2).Range.Text = "27+58="
$t.Cell(5,3).Range.Text = "93-88="
$t.Cell(5,4).Range.Text = "59+6="
$t.Cell(5,5).Range.Text = "14+57="
$t.Cell(6,1).Range.Text = "56+9="
$t.Cell(6,2).Range.Text = "47+16="
$t.Cell(6,3).Range.Text = "17+49="
$t.Cell(6,4).Range.Text = "76-69="
$t.Cell(6,5).Range.Text = "67+19="
$t.Cell(7,1).Range.Text = "78+19="
$t.Cell(7,2).Range.Text = "6+36="
$t.Cell(7,3).Range.Text = "97-8="
$t.Cell(7,4).Range.Text = "90-51="
$t.Cell(7,5).Range.Text = "29+29="
$t.Cell(8,1).Range.Text = "60-42="
$t.Cell(8,2).Range.Text = "37-18="
$t.Cell(8,3).Range.Text = "84-9="
$t.Cell(8,4).Range.Text = "52-28="
$t.Cell(8,5).Range.Text = "50-11="
$t.Cell(9,1).Range.Text = "45+19="
$t.Cell(9,2).Range.Text = "19+39="
$t.Cell(9,3).Range.Text = "47+44="
$t.Cell(9,4).Range.Text = "96-27="
$t.Cell(9,5).Range.Text = "82-4="
$t.Cell(10,1).Range.Text = "39+27="
$t.Cell(10,2).Range.Text = "33-7="
$t.Cell(10,3).Range.Text = "93-37="
$t.Cell(10,4).Range.Text = "95-49="
$t.Cell(10,5).Range.Text = "38+37="
$t.Cell(11,1).Range.Text = "38+17="
$t.Cell(11,2).Range.Text = "48+16="
$t.Cell(11,3).Range.Text = "83-35="
$t.Cell(11,4).Range.Text = "14-5="
$t.Cell(11,5).Range.Text = "46+39="
$t.Cell(12,1).Range.Text = "71-16="
$t.Cell(12,2).Range.Text = "45-8="
$t.Cell(12,3).Range.Text = "16+5="
$t.Cell(12,4).Range.Text = "40-12="
$t.Cell(12,5).Range.Text = "51-9="
$t.Cell(13,1).Range.Text = "91-52="
$t.Cell(13,2).Range.Text = "61-9="
$t.Cell(13,3).Range.Text = "54-8="
$t.Cell(13,4).Range.Text = "78-19="
$t.Cell(13,5).Range.Text = "29+24="
$t.Cell(14,1).Range.Text = "59+3="
$t.Cell(14,2).Range.Text = "29+9="
$t.Cell(14,3).Range.Text = "4+38="
$t.Cell(14,4).Range.Text = "19+7="
$t.Cell(14,5).Range.Text = "77-8="
$t.Cell(15,1).Range.Text = "83-34="
$t.Cell(15,2).Range.Text = "9+15="
$t.Cell(15,3).Range.Text = "48+9="
$t.Cell(15,4).Range.Text = "72-14="
$t.Cell(15,5).Range.Text = "19+24="
$t.Cell(16,1).Range.Text = "14+68="
$t.Cell(16,2).Range.Text = "90-67="
$t.Cell(16,3).Range.Text = "28+57="
$t.Cell(16,4).Range.Text = "37+44="
$t.Cell(16,5).Range.Text = "21-7="
$t.Cell(17,1).Range.Text = "82-27="
$t.Cell(17,2).Range.Text = "37+5="
$t.Cell(17,3).Range.Text = "37+37="
$t.Cell(17,4).Range.Text = "71-8="
$t.Cell(17,5).Range.Text = "46-17="
$t.Cell(18,1).Range.Text = "49+47="
$t.Cell(18,2).Range.Text = "73-44="
$t.Cell(18,3).Range.Text = "7+48="
$t.Cell(18,4).Range.Text = "47-39="
$t.Cell(18,5).Range.Text = "72-43="
$t.Cell(19,1).Range.Text = "83-35="
$t.Cell(19,2).Range.Text = "69+12="
$t.Cell(19,3).Range.Text = "44-36="
$t.Cell(19,4).Range.Text = "25+68="
$t.Cell(19,5).Range.Text = "16+36="
$t.Cell(20,1).Range.Text = "86-57="
$t.Cell(20,2).Range.Text = "54-7="
$t.Cell(20,3).Range.Text = "49+16="
$t.Cell(20,4).Range.Text = "23+69="
$t.Cell(20,5).Range.Text = "91-89="

Write-Host "Updated cells:" 100
